function Set-TextValue($ws, $cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

Set-TextValue $ws "D2" "51.931.57"
Set-TextValue $ws "E2" "  -0.46%  "
Set-TextValue $ws "D3" "2.790.85"
Set-TextValue $ws "E3" "  -2.00%  "
Set-TextValue $ws "E4" "  +0.00%  "
Set-TextValue $ws "D5" "361.66"
Set-TextValue $ws "E5" "  +0.02%  "
Set-TextValue $ws "D6" "109.42"
Set-TextValue $ws "E6" "  -3.80%  "
Set-TextValue $ws "E7" "  -2.31%  "
Set-TextValue $ws "E8" "  +0.05%  "
Set-TextValue $ws "E9" "  -1.89%  "
Set-TextValue $ws "D10" "40.11"
Set-TextValue $ws "E10" "  -3.60%  "
Set-TextValue $ws "D11" "0.0848"
Set-TextValue $ws "E11" "  -1.61%  "
Set-TextValue $ws "E12" "  +1.07%  "
Set-TextValue $ws "D13" "19.48"
Set-TextValue $ws "E13" "  -2.65%  "
Set-TextValue $ws "D14" "7.56"
Set-TextValue $ws "E14" "  -2.84%  "
Set-TextValue $ws "D15" "3.227.00"
Set-TextValue $ws "E15" "  -2.17%  "
Set-TextValue $ws "D16" "2.797.07"
Set-TextValue $ws "E16" "  -2.08%  "
Set-TextValue $ws "E17" "  +3.71%  "
Set-TextValue $ws "D18" "51.913.63"
Set-TextValue $ws "E18" "  -0.19%  "
Set-TextValue $ws "D19" "7.47"
Set-TextValue $ws "E19" "  -1.22%  "
Set-TextValue $ws "D20" "3.09"
Set-TextValue $ws "E20" "  -2.24%  "
Set-TextValue $ws "D21" "13.11"
Set-TextValue $ws "E21" "  -3.37%  "
$val_D22 = "{0}{1}{2}" -f "0.0", $sub3, "0976"
Set-TextValue $ws "D22" $val_D22
Set-TextValue $ws "E22" "  -1.69%  "
Set-TextValue $ws "D23" "70.35"
Set-TextValue $ws "D24" "269.82"
Set-TextValue $ws "E24" "  +0.60%  "
Set-TextValue $ws "D25" "2.76"
Set-TextValue $ws "E25" "  -2.16%  "
Set-TextValue $ws "D26" "26.54"
Set-TextValue $ws "E26" "  -2.42%  "
Set-TextValue $ws "E27" "  +0.02%  "
Set-TextValue $ws "E28" "  +15.67%  "
Set-TextValue $ws "D29" "10.27"
Set-TextValue $ws "E29" "  -1.40%  "
Set-TextValue $ws "D30" "2.27"
Set-TextValue $ws "E30" "  +0.81%  "
Set-TextValue $ws "D31" "0.0470"
Set-TextValue $ws "E31" "  +1.49%  "
Set-TextValue $ws "D32" "52.00"
Set-TextValue $ws "E32" "  -3.10%  "
Set-TextValue $ws "D33" "34.22"
Set-TextValue $ws "E33" "  +0.37%  "
Set-TextValue $ws "D34" "5.74"
Set-TextValue $ws "E34" "  -2.65%  "
Set-TextValue $ws "D35" "0.0845"
Set-TextValue $ws "E35" "  +0.05%  "
Set-TextValue $ws "D36" "5.24"
Set-TextValue $ws "E36" "  -2.70%  "
Set-TextValue $ws "E37" "  -0.03%  "
Set-TextValue $ws "D38" "19.02"
Set-TextValue $ws "E38" "  +3.65%  "
Set-TextValue $ws "D39" "3.21"
Set-TextValue $ws "E39" "  -2.28%  "
Set-TextValue $ws "D40" "1.99"
Set-TextValue $ws "E40" "  -4.06%  "
Set-TextValue $ws "D41" "2.60"
Set-TextValue $ws "E41" "  +1.33%  "
Set-TextValue $ws "E42" "  -1.87%  "
Set-TextValue $ws "E43" "  -1.04%  "
Set-TextValue $ws "B44" "Monero"
Set-TextValue $ws "C44" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D44" "119.74"
Set-TextValue $ws "E44" "  -6.62%  "
Set-TextValue $ws "B45" "EnergySwap"
Set-TextValue $ws "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D45" "22.00"
Set-TextValue $ws "E45" "  -7.91%  "
Set-TextValue $ws "D46" "2.083.60"
Set-TextValue $ws "E46" "  -1.44%  "
Set-TextValue $ws "D47" "3.25"
Set-TextValue $ws "E47" "  -4.10%  "
Set-TextValue $ws "D49" "5.79"
Set-TextValue $ws "E49" "  -1.25%  "
Set-TextValue $ws "D50" "0.950"
Set-TextValue $ws "E50" "  -5.20%  "
Set-TextValue $ws "D51" "8.86"
Set-TextValue $ws "E51" "  -2.28%  "
